$d = $word.ActiveDocument

# Replace the old device-id string (currently split across two runs and
# wrapped in a bookmark) with the new single run of text.
$d.Content.Find.Execute(
    "flutter run --release --device-id 00008030-0006091C1140802E",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "flutter run --release --device-id 00008110-001045C4018A801E", 2
)

# Remove the now-orphaned bookmark, if it survived the replace.
if ($d.Bookmarks.Exists("OLE_LINK1")) {
    $d.Bookmarks("OLE_LINK1").Delete()
}
